# Update "想去人数" (want-to-go count) values in 展览 and 全部类型 sheets,
# reflecting a refreshed data pull (gh-pages output regenerated).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 101
$ws1.Range("F4").Value = 1481
$ws1.Range("F5").Value = 195
$ws1.Range("F6").Value = 44
$ws1.Range("F7").Value = 37
$ws1.Range("F8").Value = 9806
$ws1.Range("F9").Value = 167
$ws1.Range("F10").Value = 116
$ws1.Range("F13").Value = 370
$ws1.Range("F14").Value = 6774
$ws1.Range("F15").Value = 1082
$ws1.Range("F16").Value = 625
$ws1.Range("F18").Value = 193

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 549

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 101
$ws4.Range("F4").Value = 1482
$ws4.Range("F5").Value = 195
$ws4.Range("F7").Value = 44
$ws4.Range("F8").Value = 37
$ws4.Range("F9").Value = 549
$ws4.Range("F11").Value = 9807
$ws4.Range("F12").Value = 167
$ws4.Range("F13").Value = 116
$ws4.Range("F16").Value = 370
$ws4.Range("F17").Value = 6774
$ws4.Range("F18").Value = 1082
$ws4.Range("F19").Value = 625
$ws4.Range("F21").Value = 193
